# Retrospectiva Grupo A - apply author's edits via Word COM-interop
#
# 1) "María" -> "María Taborda" (new team member added to the list)
# 2) Expand the "Lograr que las pruebas ..." bullet with extra reasoning
# 3) Expand the "Seguir mejorando la comunicación ..." bullet with extra reasoning

$d = $word.ActiveDocument

# --- 1. Add "Taborda" as María's last name -------------------------------
$d.Content.Find.Execute(
    "María",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "María Taborda",
    2
) | Out-Null

# --- 2. Flesh out the "pruebas" retro point -------------------------------
$d.Content.Find.Execute(
    "unitarias,  funcionales",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "unitarias como funcionales",
    2
) | Out-Null

$d.Content.Find.Execute(
    "la revisión y la retrospectiva del sprint. ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "la revisión y retrospectiva del sprint, para que de esa manera el desarrollo de las mismas no utilice horas adicionales que no se contemplan como parte de las horas de trabajo asignadas para cada miembro dentro del sprint. ",
    2
) | Out-Null

# --- 3. Flesh out the "comunicación" retro point --------------------------
$d.Content.Find.Execute(
    "Seguir mejorando la comunicación en el equipo. ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Seguir mejorando la comunicación en el equipo de manera que se facilite el trabajo entre los compañeros y para que el desarrollo de la aplicación se dé incremental e integralmente.",
    2
) | Out-Null
